# Rewards.xlsx edit: add NOT_GIVEN_BIRTH state group (RAND/HUNT/MULTIPLY/FLEE)
# as four new rows (18-21), plus a new "n" indicator column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows 18-21: NOT_GIVEN_BIRTH state variants -------------------------
$ws.Range("A18").Value = "NOT_GIVEN_BIRTHRAND"
$ws.Range("B18").Value = 0
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 10
$ws.Range("E18").Value = 0

$ws.Range("A19").Value = "NOT_GIVEN_BIRTHHUNT"
$ws.Range("B19").Value = 0
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 20
$ws.Range("E19").Value = 0

$ws.Range("A20").Value = "NOT_GIVEN_BIRTHMULTIPLY"
$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 0
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 20

$ws.Range("A21").Value = "NOT_GIVEN_BIRTHFLEE"
$ws.Range("B21").Value = 0
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0

# Match the numeric display style (integer format) already used by B2:E17
$ws.Range("B18:E21").NumberFormat = "0"

# --- New column F: "NOT_GIVEN_BIRTH" flag, marked "n" for every data row ---
$ws.Range("F1").Value = "NOT_GIVEN_BIRTH"
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 6).Value = "n"
}

$ws.Columns.Item(6).ColumnWidth = 17.5

# --- Selection, matching the saved view state -------------------------------
$ws.Range("H5:I5").Select()
